$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 573

$guids = @(
  "000572",
  "000573",
  "000573",
  "000574",
  "000575",
  "000576",
  "000577",
  "000578",
  "000579",
  "000580",
  "000581",
  "000582",
  "000583",
  "000584",
  "000585",
  "000586",
  "000587",
  "000588",
  "000589",
  "000590",
  "000591",
  "000592",
  "000593",
  "000594",
  "000595",
  "000596",
  "000597",
  "000598"
)

$details = @(
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 16:29:24",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 16:31:04",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 16:31:51",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 16:34:48",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 16:35:53",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 16:36:19",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:44:45",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:45:09",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:46:04",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:49:36",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:49:49",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:50:16",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:50:42",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:51:13",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:51:32",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:52:51",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 17:54:50",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 18:00:21",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 09-Jun-2023 18:01:09",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 10-Jun-2023 09:17:29",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 13:43:15",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Live test data from raw imu reading and raw leap reading. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 13:44:00",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 41B. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 13:47:27",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 32B. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 13:47:50",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 32B. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 13:48:14",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 22B. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 13:49:27",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 42E. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 13:50:14",
  "Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 42E. Script used: Interpret_IMU_And_LeapDevice_And_Motive.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 13:50:43"
)

for ($i = 0; $i -lt $guids.Length; $i++) {
  $r = $startRow + $i
  $cellA = $ws.Cells.Item($r, 1)
  $cellA.NumberFormat = "@"
  $cellA.Value = $guids[$i]
  $cellA.Style = "Normal"

  $cellB = $ws.Cells.Item($r, 2)
  $cellB.NumberFormat = "@"
  $cellB.Value = $details[$i]
  $cellB.Style = "Normal"
}

Write-Output ($ws.UsedRange.Address())
